$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-02-06 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-07 Friday", 2) | Out-Null

# Update each answer cell in the table (row-major order, matching the source data grid)
$t = $d.Tables.Item(1)
$answers = @(
    @("31+52=83","13+19=32","15+39=54","92+3=95","28+24=52"),
    @("87+2=89","56-50=6","24+43=67","97-61=36","76-48=28"),
    @("26+23=49","47-43=4","97-84=13","10+63=73","46-16=30"),
    @("73+19=92","50-27=23","9+58=67","41+10=51","79+17=96"),
    @("14+7=21","61+32=93","44+38=82","25+38=63","99-9=90"),
    @("27+48=75","68-6=62","47-26=21","73-67=6","74-56=18"),
    @("82-80=2","26+2=28","12+10=22","86-23=63","41-19=22"),
    @("55-10=45","75+23=98","33+46=79","9+40=49","65+12=77"),
    @("25-18=7","10+34=44","98-18=80","89-85=4","29+69=98"),
    @("9+36=45","72-27=45","88-70=18","6-3=3","6+52=58"),
    @("58+7=65","72+9=81","32+11=43","64+30=94","27-13=14"),
    @("91-16=75","11+18=29","78-60=18","10+24=34","12+41=53"),
    @("28-7=21","26+51=77","45+46=91","59+12=71","9+21=30"),
    @("93+1=94","19+7=26","53-35=18","47+33=80","86-57=29"),
    @("50+21=71","97-44=53","42-29=13","34+24=58","52-2=50"),
    @("39+28=67","74-10=64","21+57=78","24+16=40","27+17=44"),
    @("82+1=83","63-57=6","40-25=15","22+67=89","12+14=26"),
    @("0+90=90","91-31=60","16+53=69","45-18=27","97-52=45"),
    @("9+43=52","93-29=64","21+6=27","22-7=15","58-24=34"),
    @("89-70=19","52+40=92","26-12=14","68-3=65","83-5=78")
)

for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $answers[$r - 1][$c - 1]
    }
}

Write-Host "Done updating" ($d.Tables.Item(1).Rows.Count) "rows"
